# Apply the edits described in the commit:
#  - "collection point" sheet: add a new data row (row 2) describing a
#    Seabed collection point with dry-mate/wet-mate electrical interfaces.
#  - "dynamic cable" sheet: drop the leading blank column (the header row
#    shifts left by one and gains a proper "id [-]" header in column A).
#  - refresh the selections/viewports that Excel recorded for both sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "dynamic cable" sheet - remove the empty leading column A
# ---------------------------------------------------------------------
$wsCable = $wb.Worksheets.Item("dynamic cable")
$wsCable.Columns("A").Delete()

# ---------------------------------------------------------------------
# 2. "collection point" sheet - add new row describing a seabed
#    collection point connected through a dry-mate / wet-mate interface
# ---------------------------------------------------------------------
$wsPoint = $wb.Worksheets.Item("collection point")

$wsPoint.Range("A2").Value = 1
$wsPoint.Range("K2").Value = "dry-mate"
$wsPoint.Range("L2").Value = 4
$wsPoint.Range("N2").Value = "wet-mate"
$wsPoint.Range("O2").Value = 1
$wsPoint.Range("B2").Value = "Seabed"

$wsPoint.Range("A2").HorizontalAlignment = -4108
$wsPoint.Range("B2").HorizontalAlignment = -4108
$wsPoint.Range("K2").HorizontalAlignment = -4108
$wsPoint.Range("L2").HorizontalAlignment = -4108
$wsPoint.Range("N2").HorizontalAlignment = -4108
$wsPoint.Range("O2").HorizontalAlignment = -4108

# ---------------------------------------------------------------------
# 3. Restore/update the view state (active cell & scroll position) for
#    both worksheets, finishing on "collection point" so it stays the
#    selected tab, like in the source workbook.
# ---------------------------------------------------------------------
$wsCable.Activate()
$wsCable.Range("V8").Select() | Out-Null

$wsPoint.Activate()
$wsPoint.Range("A1").Select() | Out-Null
$wsPoint.Range("K9").Select() | Out-Null
